$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("工作表1")
$ws2old = $wb.Worksheets.Item("工作表1 (2)")

# ---------------------------------------------------------------------------
# Sheet "工作表1 (2)" (sheet2): the re-uploaded workbook shows this sheet's
# internal sheetId bumped from 11 to 12 even though its name/position/data
# are (almost) unchanged -- that only happens when the tab is recreated.
# Reproduce it by duplicating the sheet, dropping the original, and renaming
# the duplicate back to the original name, then apply the real content
# changes: the header row's C/D labels move from b/c to c/d, and the blank
# B2:D4 block gets filled with the live B$1&$A2-style concatenation formula.
# ---------------------------------------------------------------------------
$ws2idx = $ws2old.Index
$ws2old.Copy($null, $ws2old)
$wsTemp = $wb.Worksheets.Item($ws2idx + 1)
$ws2old.Name = "TEMP_工作表1_2_OLD"
$wsTemp.Name = "工作表1 (2)"
$ws2old.Delete()
$ws2 = $wb.Worksheets.Item("工作表1 (2)")

$ws2.Range("C1").Value = "c"
$ws2.Range("D1").Value = "d"

$ws2.Range("C2").Formula = '=C$1&$A2'
$ws2.Range("B2").Formula = '=B$1&$A2'
$ws2.Range("D2").Formula = '=D$1&$A2'
$ws2.Range("B3:D4").Formula = '=B$1&$A3'

$ws2.Range("D8").Select()

# ---------------------------------------------------------------------------
# Sheet "工作表1" (sheet1): drop the old A1 header/"#" column, rename the
# B/C/D headers to b/c/d, renumber the A-column counters 2..4, drop the old
# row 5, and replace the literal "updated" fills in B2:D4 with the live
# concatenation formula B$1&$A2 filled down. Edited last so it ends up the
# active sheet/tab again, matching the original tabSelected state.
# ---------------------------------------------------------------------------
$ws1.Activate()

$ws1.Range("A1").ClearContents()
$ws1.Range("B1").Value = "b"
$ws1.Range("C1").Value = "c"
$ws1.Range("D1").Value = "d"

$ws1.Range("A2").Value = 2
$ws1.Range("A3").Value = 3
$ws1.Range("A4").Value = 4

$ws1.Range("B2:D4").Formula = '=B$1&$A2'

$ws1.Range("A5:D5").ClearContents()

$ws1.Range("G10").Select()
